{"js": "const replacements = [\n  [\"2026-01-11 Sunday\", \"2026-01-12 Monday\"],\n  [\"510\u00f73=\", \"681\u00f77=\"],\n  [\"338\u00f77=\", \"435\u00f72=\"],\n  [\"280\u00f73=\", \"154\u00f74=\"],\n  [\"876\u00f77=\", \"895\u00f76=\"],\n  [\"543\u00f79=\", \"565\u00f78=\"],\n  [\"269\u00f72=\", \"714\u00f72=\"],\n  [\"650\u00f78=\", \"720\u00f79=\"],\n  [\"629\u00f75=\", \"279\u00f77=\"],\n  [\"255\u00f73=\", \"222\u00f74=\"],\n  [\"344\u00f73=\", \"531\u00f77=\"],\n  [\"231\u00f74=\", \"980\u00f77=\"],\n  [\"941\u00f76=\", \"701\u00f72=\"],\n  [\"618\u00f79=\", \"349\u00f76=\"],\n  [\"497\u00f74=\", \"145\u00f77=\"],\n  [\"169\u00f78=\", \"508\u00f74=\"],\n  [\"838\u00f73=\", \"723\u00f77=\"],\n  [\"605\u00f78=\", \"582\u00f75=\"],\n  [\"378\u00f78=\", \"286\u00f78=\"],\n  [\"661\u00f73=\", \"332\u00f79=\"],\n  [\"499\u00f73=\", \"812\u00f79=\"],\n  [\"882\u00f74=\", \"773\u00f79=\"],\n  [\"468\u00f73=\", \"603\u00f78=\"],\n  [\"171\u00f78=\", \"293\u00f73=\"],\n  [\"385\u00f76=\", \"751\u00f79=\"],\n  [\"419\u00f76=\", \"409\u00f77=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2026-01-11 Sunday\", \"2026-01-12 Monday\"),\n    @(\"510\u00f73=\", \"681\u00f77=\"),\n    @(\"338\u00f77=\", \"435\u00f72=\"),\n    @(\"280\u00f73=\", \"154\u00f74=\"),\n    @(\"876\u00f77=\", \"895\u00f76=\"),\n    @(\"543\u00f79=\", \"565\u00f78=\"),\n    @(\"269\u00f72=\", \"714\u00f72=\"),\n    @(\"650\u00f78=\", \"720\u00f79=\"),\n    @(\"629\u00f75=\", \"279\u00f77=\"),\n    @(\"255\u00f73=\", \"222\u00f74=\"),\n    @(\"344\u00f73=\", \"531\u00f77=\"),\n    @(\"231\u00f74=\", \"980\u00f77=\"),\n    @(\"941\u00f76=\", \"701\u00f72=\"),\n    @(\"618\u00f79=\", \"349\u00f76=\"),\n    @(\"497\u00f74=\", \"145\u00f77=\"),\n    @(\"169\u00f78=\", \"508\u00f74=\"),\n    @(\"838\u00f73=\", \"723\u00f77=\"),\n    @(\"605\u00f78=\", \"582\u00f75=\"),\n    @(\"378\u00f78=\", \"286\u00f78=\"),\n    @(\"661\u00f73=\", \"332\u00f79=\"),\n    @(\"499\u00f73=\", \"812\u00f79=\"),\n    @(\"882\u00f74=\", \"773\u00f79=\"),\n    @(\"468\u00f73=\", \"603\u00f78=\"),\n    @(\"171\u00f78=\", \"293\u00f73=\"),\n    @(\"385\u00f76=\", \"751\u00f79=\"),\n    @(\"419\u00f76=\", \"409\u00f77=\")\n)\n\nforeach ($pair in $replacements) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
